# Bereich 'Jobtitel' in Datenbank inkl. entsprechende Stored Procedures erstellt.
# Datenimport und alle Tests laufen.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---------------------------------------------------------------------------
# 1) Mitarbeitertyp (row 23): "Werkstudent" -> "Angestellter"
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = "Angestellter"

# ---------------------------------------------------------------------------
# 2) Rows 26-28 (Abteilung / Abteilungsabkürzung / Führungskraft): the "filled
#    out" indicator in column D switches from the red block to the green one.
#    Column D carries no value, only fill formatting, so copy that formatting
#    from a cell that already has the green fill (e.g. D25).
# ---------------------------------------------------------------------------
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D26:D28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------------
# 3) Insert the new "Jobtitel" / "Erfahrungsstufe" pair just above the
#    existing "Gesellschaft" row, pushing "Gesellschaft" one row down.
#    Row 29 becomes "Jobtitel" (new answer "Data Analyst"),
#    row 30 becomes "Erfahrungsstufe" (new answer "Junior"),
#    row 31 becomes "Gesellschaft" (still unanswered).
# ---------------------------------------------------------------------------

# Row 29: take on the "orange" look that row 30/31 ("Jobtitel") used to have.
$ws.Range("A30:B30").Copy() | Out-Null
$ws.Range("A29:B29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A29").Value = "Jobtitel"
$ws.Range("B29").Value = "Data Analyst"

# Row 30 keeps its existing formatting, only the label/value change.
$ws.Range("A30").Value = "Erfahrungsstufe"
$ws.Range("B30").Value = "Junior"

# Columns D29/D30 flip from the red block to the green one, same as D26:D28.
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D29:D30").PasteSpecial(-4122) | Out-Null

# Row 31: now holds "Gesellschaft" again, with the formatting that row 29
# used to have (still empty, thick-bottom border stays, since that is a row
# level property we never touched).
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A31:B31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A31").Value = "Gesellschaft"
$ws.Range("B31").Value = $null

# ---------------------------------------------------------------------------
# 4) Append a brand-new row 48 "Anzahl Kinder" after the existing last row.
# ---------------------------------------------------------------------------
$ws.Range("A47:B47").Copy() | Out-Null
$ws.Range("A48:B48").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D47").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4122) | Out-Null

$ws.Range("A48").Value = "Anzahl Kinder"
$ws.Range("B48").Value = $null

# ---------------------------------------------------------------------------
# 5) Update the view: scroll so row 16 is at the top and select C39.
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("C39").Select()
